$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 280
$ws1.Range("F4").Value = 1270
$ws1.Range("F7").Value = 61
$ws1.Range("F9").Value = 147
$ws1.Range("F10").Value = 3507
$ws1.Range("F11").Value = 135
$ws1.Range("F12").Value = 90
$ws1.Range("F14").Value = 45
$ws1.Range("F16").Value = 607
$ws1.Range("F18").Value = 757
$ws1.Range("F20").Value = 123
$ws1.Range("F23").Value = 68
$ws1.Range("F24").Value = 2682
$ws1.Range("F25").Value = 5186
$ws1.Range("F29").Value = 3077
$ws1.Range("F30").Value = 285
$ws1.Range("F35").Value = 125
$ws1.Range("F37").Value = 314
$ws1.Range("F38").Value = 29
$ws1.Range("F40").Value = 805
$ws1.Range("F42").Value = 3

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 280
$ws4.Range("F4").Value = 1270
$ws4.Range("F7").Value = 61
$ws4.Range("F9").Value = 147
$ws4.Range("F10").Value = 3507
$ws4.Range("F11").Value = 135
$ws4.Range("F12").Value = 90
$ws4.Range("F15").Value = 45
$ws4.Range("F17").Value = 607
$ws4.Range("F19").Value = 757
$ws4.Range("F21").Value = 123
$ws4.Range("F24").Value = 68
$ws4.Range("F25").Value = 2683
$ws4.Range("F26").Value = 5186
$ws4.Range("F30").Value = 3077
$ws4.Range("F31").Value = 285
$ws4.Range("F36").Value = 125
$ws4.Range("F38").Value = 314
$ws4.Range("F39").Value = 29
$ws4.Range("F41").Value = 805
$ws4.Range("F43").Value = 3
